$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append (rows 7-9)
$data = @(
    @{ Row = 7; A = 9913.9699999999993; B = 9895.17;              C = 78.05;              D = 78.2;               F = 0.19;  G = 42613.766342592593; H = $true  },
    @{ Row = 8; A = 9995.26;             B = 9913.9699999999993;  C = 77.739999999999995; D = 78.38;              F = 0.82;  G = 42614.674398148149; H = $true  },
    @{ Row = 9; A = 9941.2900000000009;  B = 9995.26;             C = 78.36;              D = 77.94;              F = -0.54; G = 42615.752523148149; H = $false }
)

foreach ($item in $data) {
    $r = $item.Row

    # Match the date/time style already used by the existing G column cells
    $ws.Cells.Item(6, 7).Copy() | Out-Null
    $ws.Cells.Item($r, 7).PasteSpecial(-4122) | Out-Null # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $false
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
